$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header I1 ("instax" -> "Solo/Group")
$ws.Range("I1").Value = "Solo/Group"

# Add new row 3 with data
$ws.Range("A3").Value = "DONE CHOOSING"
$ws.Range("B3").Value = "ASD"
$ws.Range("C3").Value = "asd@asd.asd"
$ws.Range("D3").Value = "Basic"
$ws.Range("E3").Value = "ABOGADO, MAISA C.jpeg"
$ws.Range("F3").Value = "N/A"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "N/A"
$ws.Range("I3").Value = "N/A"
$ws.Range("J3").Value = "2026-01-16 00:51:24"

# Match row 2's formatting (thin border all around) on the new row
$rng = $ws.Range("A3:J3")
$rng.Borders.Item(1).LineStyle = 1
$rng.Borders.Item(2).LineStyle = 1
$rng.Borders.Item(3).LineStyle = 1
$rng.Borders.Item(4).LineStyle = 1

# Extend the AutoFilter range from A1:J2 to A1:J3.
# Toggle off first (re-applying the same range flips AutoFilter off),
# then re-apply with the new, larger range.
$ws.Range("A1:J2").AutoFilter()
$ws.Range("A1:J3").AutoFilter()
